$d = $word.ActiveDocument

# Locate the "LOB1036: Geometria Analítica (Requisito fraco)" paragraph
# (the last line of the "Requisitos" section that must be kept).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "LOB1036: Geometria Analítica (Requisito fraco)") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    # The four paragraphs right after it are removed:
    #   - an empty paragraph
    #   - "Ver no Jupiter Salvar em pdf Salvar em docx"
    #   - an empty paragraph
    #   - an empty, page-break-before paragraph
    # leaving the following empty / page-break paragraphs intact.
    $startPara = $d.Paragraphs.Item($targetIndex + 1)
    $endPara = $d.Paragraphs.Item($targetIndex + 4)

    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRange.Delete()
}
